# Resize the two ASCII-art "box" borders (symmetrical segments) and
# update the two "First Segment" derivation numbers (PBC / existential
# elimination fix) in the natural-deduction table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 4  -> box #1 top border  (┌---…---┐)
$t.Rows.Item(4).Cells.Item(1).Range.Text = `
  "                  ┌--------------------------------------------------------------------------------------------------------------┐"

# Row 8  -> box #2 top border  (┌---…---┐)
$t.Rows.Item(8).Cells.Item(1).Range.Text = `
  "                        ┌------------------------------------------------------------------------------------------------------┐"

# Row 12 -> box #2 bottom border (└---…---┘)
$t.Rows.Item(12).Cells.Item(1).Range.Text = `
  "                        └------------------------------------------------------------------------------------------------------┘"

# Row 13 -> line 9, "First Segment" column: 6 -> 5
$t.Rows.Item(13).Cells.Item(4).Range.Text = "5"

# Row 14 -> box #1 bottom border (└---…---┘)
$t.Rows.Item(14).Cells.Item(1).Range.Text = `
  "                  └--------------------------------------------------------------------------------------------------------------┘"

# Row 15 -> line 10, "First Segment" column: 1 -> 2
$t.Rows.Item(15).Cells.Item(4).Range.Text = "2"
